# chore: update Sheets via scheduled runner
# Refreshes market-board derived columns (H:N) on the Leve profit sheets
# with the latest pricing snapshot. Values only - no formulas/formatting.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 397.05554
$ws.Range("I9").Value = 314.07144
$ws.Range("K9").Value = 314.07144
$ws.Range("M9").Value = -145.07144

$ws.Range("H41").Value = 7814389
$ws.Range("I41").Value = 11364979
$ws.Range("J41").Value = 3090.2
$ws.Range("K41").Value = 11364979
$ws.Range("L41").Value = 3090.2
$ws.Range("M41").Value = -11364539
$ws.Range("N41").Value = -3970.2

$ws.Range("H62").Value = 1670
$ws.Range("I62").Value = 1505
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1505
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -881
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 1670
$ws.Range("I65").Value = 1505
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 7525
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -4405
$ws.Range("N65").Value = -16240

$ws.Range("H138").Value = 2994.51
$ws.Range("I138").Value = 1194.8363
$ws.Range("J138").Value = 5194.1113
$ws.Range("K138").Value = 3584.5089
$ws.Range("L138").Value = 15582.3339
$ws.Range("M138").Value = 1555.4911
$ws.Range("N138").Value = -25862.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 51713.145
$ws.Range("I74").Value = 113787.11
$ws.Range("K74").Value = 113787.11
$ws.Range("M74").Value = -112913.11

$ws.Range("H77").Value = 51713.145
$ws.Range("I77").Value = 113787.11
$ws.Range("K77").Value = 568935.55
$ws.Range("M77").Value = -564567.55

$ws.Range("H122").Value = 3936.8572
$ws.Range("I122").Value = 3358.8572
$ws.Range("J122").Value = 4514.857
$ws.Range("K122").Value = 10076.5716
$ws.Range("L122").Value = 13544.571
$ws.Range("M122").Value = -7626.571599999999
$ws.Range("N122").Value = -18444.571

$ws.Range("H132").Value = 3999.9285
$ws.Range("I132").Value = 1544.0256
$ws.Range("K132").Value = 4632.0768
$ws.Range("M132").Value = -2102.0768

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H86").Value = 8966691
$ws.Range("J86").Value = 3434.25
$ws.Range("L86").Value = 3434.25
$ws.Range("N86").Value = -5680.25

$ws.Range("H89").Value = 8966691
$ws.Range("J89").Value = 3434.25
$ws.Range("L89").Value = 17171.25
$ws.Range("N89").Value = -28403.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3646.1667
$ws.Range("I99").Value = 2667.0527
$ws.Range("K99").Value = 2667.0527
$ws.Range("M99").Value = -1169.0527

$ws.Range("H126").Value = 3646.1667
$ws.Range("I126").Value = 2667.0527
$ws.Range("K126").Value = 8001.158100000001
$ws.Range("M126").Value = -5531.158100000001

$ws.Range("H132").Value = 7493.5557
$ws.Range("I132").Value = 4486.9165
$ws.Range("K132").Value = 13460.7495
$ws.Range("M132").Value = -10930.7495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 953.2
$ws.Range("I70").Value = 953.2
$ws.Range("K70").Value = 2859.6
$ws.Range("M70").Value = -2544.6

$ws.Range("H73").Value = 953.2
$ws.Range("I73").Value = 953.2
$ws.Range("K73").Value = 2859.6
$ws.Range("M73").Value = -1767.6

$ws.Range("H75").Value = 35095300
$ws.Range("I75").Value = 47620970
$ws.Range("J75").Value = 27788660
$ws.Range("K75").Value = 142862910
$ws.Range("L75").Value = 83365980
$ws.Range("M75").Value = -142861912
$ws.Range("N75").Value = -83367976

$ws.Range("H78").Value = 35095300
$ws.Range("I78").Value = 47620970
$ws.Range("J78").Value = 27788660
$ws.Range("K78").Value = 428588730
$ws.Range("L78").Value = 250097940
$ws.Range("M78").Value = -428583738
$ws.Range("N78").Value = -250107924

$ws.Range("H80").Value = 35719010
$ws.Range("I80").Value = 27781624
$ws.Range("K80").Value = 83344872
$ws.Range("M80").Value = -83343936

$ws.Range("H83").Value = 35719010
$ws.Range("I83").Value = 27781624
$ws.Range("K83").Value = 250034616
$ws.Range("M83").Value = -250029936

$ws.Range("H103").Value = 717
$ws.Range("I103").Value = 233.83333
$ws.Range("J103").Value = 1296.8
$ws.Range("K103").Value = 701.49999
$ws.Range("L103").Value = 3890.4
$ws.Range("M103").Value = 177.50001
$ws.Range("N103").Value = -5648.4

$ws.Range("H131").Value = 49213.332
$ws.Range("J131").Value = 54288.473
$ws.Range("L131").Value = 162865.419
$ws.Range("N131").Value = -172945.419

$ws.Range("H140").Value = 128188.94
$ws.Range("I140").Value = 155346.39
$ws.Range("J140").Value = 10506.667
$ws.Range("K140").Value = 466039.17
$ws.Range("L140").Value = 31520.001
$ws.Range("M140").Value = -460859.17
$ws.Range("N140").Value = -41880.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4091.5
$ws.Range("I102").Value = 3920.4119
$ws.Range("K102").Value = 3920.4119
$ws.Range("M102").Value = -2298.4119

$ws.Range("H122").Value = 8050505
$ws.Range("I122").Value = 9056318
$ws.Range("K122").Value = 27168954
$ws.Range("M122").Value = -27166504

$ws.Range("H126").Value = 5350.3
$ws.Range("I126").Value = 2303.3333
$ws.Range("K126").Value = 6909.999899999999
$ws.Range("M126").Value = -4439.999899999999

$ws.Range("H132").Value = 3917.8965
$ws.Range("I132").Value = 1471.2222
$ws.Range("J132").Value = 7921.5454
$ws.Range("K132").Value = 4413.6666
$ws.Range("L132").Value = 23764.6362
$ws.Range("M132").Value = -1883.6666
$ws.Range("N132").Value = -28824.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4290.727
$ws.Range("J7").Value = 5075.5
$ws.Range("L7").Value = 5075.5
$ws.Range("N7").Value = -5299.5

$ws.Range("H40").Value = 5601.381
$ws.Range("I40").Value = 4636
$ws.Range("J40").Value = 6888.5557
$ws.Range("K40").Value = 4636
$ws.Range("L40").Value = 6888.5557
$ws.Range("M40").Value = -4500
$ws.Range("N40").Value = -7160.5557

$ws.Range("H46").Value = 1839.3334
$ws.Range("I46").Value = 1603.3914
$ws.Range("J46").Value = 3196
$ws.Range("K46").Value = 1603.3914
$ws.Range("L46").Value = 3196
$ws.Range("M46").Value = -1415.3914
$ws.Range("N46").Value = -3572

$ws.Range("H68").Value = 8248.833000000001
$ws.Range("J68").Value = 8600
$ws.Range("L68").Value = 8600
$ws.Range("N68").Value = -10098

$ws.Range("H71").Value = 8248.833000000001
$ws.Range("J71").Value = 8600
$ws.Range("L71").Value = 43000
$ws.Range("N71").Value = -50488

$ws.Range("H93").Value = 1274.6666
$ws.Range("I93").Value = 1274.6666
$ws.Range("K93").Value = 1274.6666
$ws.Range("M93").Value = -26.66660000000002

$ws.Range("H100").Value = 3981.3333
$ws.Range("I100").Value = 3972
$ws.Range("K100").Value = 3972
$ws.Range("M100").Value = -3431

$ws.Range("H126").Value = 4290.727
$ws.Range("J126").Value = 5075.5
$ws.Range("L126").Value = 15226.5
$ws.Range("N126").Value = -20166.5

$ws.Range("H132").Value = 9621510
$ws.Range("I132").Value = 18521890
$ws.Range("J132").Value = 9101.639999999999
$ws.Range("K132").Value = 55565670
$ws.Range("L132").Value = 27304.92
$ws.Range("M132").Value = -55563140
$ws.Range("N132").Value = -32364.92

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19253.166
$ws.Range("J41").Value = 17646.8
$ws.Range("L41").Value = 17646.8
$ws.Range("N41").Value = -18426.8

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 5676.2563
$ws.Range("I132").Value = 5441.643
$ws.Range("K132").Value = 16324.929
$ws.Range("M132").Value = -13794.929

$ws.Range("H136").Value = 21493520
$ws.Range("I136").Value = 45456830
$ws.Range("J136").Value = 405804.9
$ws.Range("K136").Value = 136370490
$ws.Range("L136").Value = 1217414.7
$ws.Range("M136").Value = -136367940
$ws.Range("N136").Value = -1222514.7
